$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AA5").Value = "test"
$ws.Range("AA5").ClearFormats()
$ws.Range("AA5").Font.Size = 11
$ws.Range("AA5").Font.ColorIndex = -4105
